# Add two new weekly scoreboard rows (162 and 163) to the bottom of the
# tracker, matching the existing data layout:
# A=Participant, B=Date, C=Workout Type, D=Total Duration, E=Total Distance,
# F=Total Elevation, G=Zone 1, H=Zone 2, I=Zone 3, J=Zone 4, K=Zone 5,
# L=Workout Level, M=Week

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 162 ---
$ws.Range("A162").Value = "Jeremiah"
$ws.Range("B162").Value = 45478
$ws.Range("C162").Value = "Workout"
$ws.Range("D162").Value = 59
$ws.Range("E162").Value = 0
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 41
$ws.Range("H162").Value = 15
$ws.Range("I162").Value = 2
$ws.Range("J162").Value = 0
$ws.Range("K162").Value = 0
$ws.Range("L162").Value = "Agile Antelope"
$ws.Range("M162").Value = 4

# --- Row 163 ---
$ws.Range("A163").Value = "Steven"
$ws.Range("B163").Value = 45478
$ws.Range("C163").Value = "Walk"
$ws.Range("D163").Value = 35
$ws.Range("E163").Value = 1.4
$ws.Range("F163").Value = 85
$ws.Range("G163").Value = 35
$ws.Range("H163").Value = 0
$ws.Range("I163").Value = 0
$ws.Range("J163").Value = 0
$ws.Range("K163").Value = 0
$ws.Range("L163").Value = "Mighty Monkey"
$ws.Range("M163").Value = 4

# Copy the date number format from the cell directly above the new rows so
# the same shared date style (not a brand-new number format) is reused.
$ws.Range("B161").Copy()
$ws.Range("B162:B163").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-home the frozen header pane on the newly added bottom rows and park
# the selection on the next empty row, same as Excel does after data entry.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A150").Select()
$ws.Range("A164").Select()
